$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @"
2,7.47086763381958,6.874066829681396,6.439445495605469,6.856260776519775,4.709311962127686,5.684148788452148
3,8.163686752319336,7.154438495635986,7.169380187988281,7.617014408111572,5.27050256729126,6.012582778930664
4,8.051456451416016,7.001540660858154,6.663507461547852,7.449251174926758,6.162642955780029,6.38532543182373
5,6.937328338623047,6.893825054168701,5.490668296813965,6.557094097137451,4.756118297576904,5.262788772583008
6,7.428953647613525,7.08682107925415,6.104562759399414,6.653888702392578,4.73939847946167,5.076407909393311
7,6.89729118347168,6.874671459197998,5.167336940765381,5.854881286621094,4.248826503753662,4.442036628723145
8,7.927778720855713,7.318855762481689,7.045222759246826,7.29020881652832,4.63810396194458,6.127086639404297
9,6.582281112670898,6.670324802398682,4.935073375701904,5.760617256164551,3.9692223072052,4.789995193481445
10,8.167893409729004,7.179366111755371,6.776599884033203,7.803391933441162,6.678578853607178,6.846604824066162
11,7.427693367004395,6.885190963745117,6.554356098175049,6.670352458953857,4.516287326812744,7.39923095703125
12,7.796838760375977,7.918933391571045,6.626201629638672,7.223772525787354,6.013787746429443,5.638765811920166
13,8.210800170898438,7.990981101989746,7.335238456726074,7.389788627624512,5.113366603851318,5.597140312194824
14,7.792764186859131,7.128057956695557,6.478930473327637,7.25791072845459,5.982015609741211,6.373027324676514
15,7.983171939849854,7.927110195159912,7.005264759063721,7.346334457397461,5.503777027130127,5.602215766906738
16,7.66686487197876,7.075197219848633,6.777567863464355,6.995403289794922,4.843403339385986,7.263015747070312
17,7.317229747772217,7.642148494720459,5.819735527038574,6.64349365234375,4.674530029296875,5.178330898284912
18,7.398324012756348,7.252796649932861,6.012967109680176,6.790703773498535,5.755001544952393,6.771960258483887
19,7.522325992584229,7.375320911407471,6.598684310913086,6.755066394805908,4.846784591674805,5.997751235961914
20,8.508449554443359,7.555001735687256,7.724048137664795,8.000906944274902,5.762430667877197,7.387633323669434
21,7.513181209564209,7.581302642822266,5.994894504547119,6.740391731262207,4.699885845184326,5.243282318115234
22,7.995851993560791,7.147819042205811,6.777787685394287,7.417803764343262,5.96380090713501,7.000676155090332
23,8.466668128967285,8.192058563232422,7.642704963684082,7.653449535369873,5.350515365600586,5.490033149719238
24,7.346930980682373,7.371703624725342,6.097207069396973,6.845163822174072,5.018463611602783,4.899500846862793
25,7.380749225616455,6.749911785125732,5.772681713104248,7.016166687011719,5.411291599273682,4.851583957672119
26,7.01768159866333,6.85221529006958,5.402104377746582,6.161467552185059,4.442221164703369,4.622466564178467
27,8.010285377502441,8.606967926025391,5.960823059082031,7.169922828674316,5.864020824432373,4.482544898986816
28,8.515973091125488,7.66918420791626,7.040867805480957,8.179985046386719,7.189237594604492,5.947932243347168
29,7.333739280700684,7.427762985229492,5.786116600036621,6.514212608337402,4.728014469146729,4.450439453125
30,8.293649673461914,8.222844123840332,6.518852233886719,7.685305118560791,6.952371120452881,4.609456539154053
31,7.392475605010986,6.733833789825439,6.252872943878174,6.85020923614502,5.104836940765381,6.140723705291748
32,8.112260818481445,7.644919872283936,7.386974811553955,7.539396286010742,5.11543083190918,6.603957176208496
33,8.863296508789062,8.425434112548828,8.246442794799805,8.301699638366699,6.409745693206787,7.656247138977051
34,6.770586490631104,6.879781723022461,5.216084957122803,6.028234481811523,4.565867900848389,4.727926254272461
35,6.930089473724365,7.406124114990234,5.231637477874756,6.274822235107422,4.847255229949951,4.552416801452637
36,7.411664485931396,7.553601264953613,6.34324836730957,6.693538665771484,4.937832355499268,5.417017936706543
37,7.450937747955322,7.15873384475708,6.401279926300049,6.900914669036865,4.580910205841064,5.984614372253418
38,8.625436782836914,8.351720809936523,8.001162528991699,7.752532005310059,5.29693603515625,5.888449668884277
39,7.806265354156494,8.275774955749512,6.2722487449646,6.837416648864746,5.455821514129639,4.515962600708008
40,8.814878463745117,8.43437385559082,8.214629173278809,8.16465950012207,6.174551963806152,6.758768558502197
41,6.669003963470459,6.71594762802124,5.334782123565674,6.193369388580322,4.706305027008057,6.100512981414795
42,7.682479858398438,6.91636323928833,6.04643726348877,7.353627681732178,6.646446704864502,5.930568695068359
43,7.190935611724854,6.67522668838501,5.96125316619873,6.882350444793701,5.072341442108154,6.575161457061768
44,8.552821159362793,8.218327522277832,7.719227313995361,7.875578880310059,6.214391231536865,5.901785373687744
45,8.300132751464844,6.948515892028809,7.077082633972168,7.836573600769043,5.924401760101318,6.171560764312744
46,7.15196704864502,6.653868198394775,5.842792510986328,6.578197956085205,5.138564109802246,6.031874656677246
47,7.94158411026001,6.795178890228271,6.94889497756958,7.283754348754883,5.262527942657471,6.739367485046387
48,8.05488395690918,7.667623996734619,7.231516361236572,7.263659954071045,5.195600032806396,5.58540678024292
49,8.398725509643555,7.748189449310303,7.360597133636475,7.918464660644531,5.986591815948486,5.709778308868408
50,8.521249771118164,7.817073345184326,7.619759559631348,8.030502319335938,6.200858592987061,7.477837085723877
51,7.826517105102539,7.637566566467285,6.635591506958008,7.328590393066406,6.148102283477783,6.250245094299316
52,7.119001388549805,7.30377721786499,5.171795845031738,6.397603511810303,5.374392032623291,4.558778762817383
53,7.919630527496338,7.562905788421631,6.291118621826172,7.260937690734863,5.566336154937744,4.640822410583496
54,7.129404544830322,6.8351731300354,5.566332817077637,6.580592155456543,4.843445301055908,5.282955646514893
55,8.226907730102539,8.484601974487305,7.660160064697266,7.578631401062012,5.451854228973389,6.367569446563721
56,6.746710300445557,6.505619049072266,5.360849857330322,6.06002950668335,4.385526180267334,4.594614028930664
57,7.646710872650146,7.000685214996338,6.620978832244873,6.929399490356445,4.705392837524414,5.491735458374023
58,8.855945587158203,8.212162971496582,8.24342155456543,8.302885055541992,6.101729869842529,6.926557540893555
59,8.64103889465332,8.034981727600098,8.134922981262207,8.177105903625488,5.751770496368408,7.409085273742676
60,7.484960079193115,7.131712436676025,5.704475402832031,7.226638793945312,5.80454683303833,5.326168060302734
61,7.242800712585449,7.304914951324463,5.532867431640625,6.564404487609863,5.498175144195557,4.480834007263184
62,7.796838760375977,7.918933391571045,6.626201629638672,7.223772525787354,6.013787746429443,5.638765811920166
63,7.979450702667236,7.402157306671143,6.882497787475586,7.546627044677734,6.053736209869385,7.229357719421387
64,7.755846977233887,7.757835865020752,6.720773220062256,7.642871856689453,5.658026695251465,6.656674385070801
65,7.259077548980713,7.363428115844727,5.885808944702148,6.549985885620117,5.366577625274658,6.608119964599609
66,7.931739330291748,8.065883636474609,7.412434577941895,7.164750576019287,4.959272861480713,7.754214763641357
67,7.156197547912598,6.66966438293457,5.086333274841309,6.293149948120117,3.973208665847778,4.380451679229736
68,8.115157127380371,7.132742404937744,7.183685302734375,7.504464149475098,5.6062331199646,6.394939422607422
69,6.95361328125,6.33339262008667,5.445131301879883,6.466211318969727,5.033066749572754,6.351366519927979
70,7.640227794647217,7.620574474334717,6.457978248596191,6.763078212738037,4.62246561050415,4.769268989562988
71,7.806265354156494,8.275774955749512,6.2722487449646,6.837416648864746,5.455821514129639,4.515962600708008
72,6.967820644378662,6.69730281829834,5.735086917877197,6.387846946716309,4.419198513031006,6.655695915222168
73,7.525012969970703,7.292644023895264,6.699174880981445,6.838784217834473,4.66582202911377,5.660482406616211
74,7.566303730010986,6.934024333953857,6.626898765563965,6.890563488006592,4.782983303070068,6.493093967437744
75,7.773093700408936,7.084635257720947,6.763698101043701,7.134882926940918,5.074550151824951,6.033588409423828
76,8.598039627075195,7.639785766601562,7.627642154693604,8.044203758239746,6.299493312835693,7.852786064147949
77,6.909545421600342,6.366891860961914,5.043941020965576,6.145024299621582,4.144981384277344,4.376417636871338
78,8.280960083007812,7.476271152496338,7.151237010955811,7.636371612548828,5.516490936279297,5.332467079162598
79,8.77407169342041,7.663079738616943,8.029826164245605,8.133111000061035,5.768377780914307,7.286948204040527
80,6.642122268676758,7.234349250793457,4.296267032623291,5.60756254196167,5.148536205291748,3.97539234161377
81,8.055285453796387,7.252442836761475,6.482804298400879,7.051577568054199,5.193886756896973,4.302562236785889
82,7.651881217956543,7.71567964553833,6.297039031982422,7.04381799697876,5.861557483673096,7.337143898010254
83,6.846999645233154,6.246854305267334,5.504982948303223,5.932002067565918,4.752709865570068,6.997433185577393
84,8.181288719177246,7.848267078399658,7.569042205810547,7.337041854858398,5.03678035736084,6.821465492248535
85,6.744333267211914,6.476171970367432,5.179632186889648,5.658768653869629,4.265886306762695,4.173860549926758
86,8.361747741699219,8.143915176391602,7.696229457855225,7.794710159301758,5.881853580474854,6.380668163299561
"@

$lines = $data -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split ","
    $row = [int]$parts[0]
    $ws.Cells.Item($row, 2).Value = [double]$parts[1]
    $ws.Cells.Item($row, 4).Value = [double]$parts[2]
    $ws.Cells.Item($row, 6).Value = [double]$parts[3]
    $ws.Cells.Item($row, 8).Value = [double]$parts[4]
    $ws.Cells.Item($row, 10).Value = [double]$parts[5]
    $ws.Cells.Item($row, 12).Value = [double]$parts[6]
}

Write-Host "Done updating $($lines.Count) rows"
